# Automatische test-sync: 2025-08-01 23:56:50
# Adds the new "Testmail #10" row to the Logs sheet, extends the
# conditional-formatting ranges to cover it, and bumps the "Overig"
# tally on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 15 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(15, 1).Value = "Is er al nieuws?"
$logs.Cells.Item(15, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(15, 3).Value = "Testmail #10: Is er al nieuws?"
$logs.Cells.Item(15, 4).Value = "Overig"
$logs.Cells.Item(15, 5).Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$logs.Cells.Item(15, 6).Value = "2025-08-01 23:56:11"
$logs.Cells.Item(15, 7).Value = "Ja"
$logs.Cells.Item(15, 8).Value = "Ja"
$logs.Cells.Item(15, 9).Value = "Nee"
$logs.Cells.Item(15, 10).Value = "Nee"

# --- Extend the conditional formatting ranges from row 14 to row 15 ------
$logs.Range("D2:D14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D15"))
$logs.Range("G2:G14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G15"))
$logs.Range("H2:H14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H15"))
$logs.Range("I2:I14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I15"))
$logs.Range("J2:J14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J15"))

# --- Dashboard sheet: bump the "Overig" count from 7 to 8 -----------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 8
